# Appends four new reference rows to the end of the single table in the
# document (the "References" table), matching the newly-added rows for
# Tkinter/Matplotlib related links.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Add-ReferenceRow($col1, $col2, $url, $date, $extraTrailingParagraph) {
    $row = $t.Rows.Add()

    $c1 = $row.Cells.Item(1)
    $c1.Range.Text = $col1

    $c2 = $row.Cells.Item(2)
    $c2.Range.Text = $col2

    # Column 3 holds a hyperlink whose visible text is the URL itself.
    # Put the URL text in first, then shrink the range by the trailing
    # paragraph mark before handing it to Hyperlinks.Add so the link wraps
    # the existing run cleanly (no stray empty run left behind).
    $c3 = $row.Cells.Item(3)
    $c3.Range.Text = $url
    $linkRange = $c3.Range
    $linkRange.MoveEnd(1, -1)
    $d.Hyperlinks.Add($linkRange, $url) | Out-Null

    $c4 = $row.Cells.Item(4)
    $c4.Range.Text = $date
    if ($extraTrailingParagraph) {
        $dateRange = $c4.Range
        $dateRange.Collapse(0)
        $dateRange.InsertParagraphAfter()
    }
}

Add-ReferenceRow `
    "How to embed Matplotlib figure in tkinter " `
    "Website that shows how to add a matplotlib figure into tkinter " `
    "https://www.geeksforgeeks.org/how-to-embed-matplotlib-charts-in-tkinter-gui/" `
    "26/4/2024" `
    $false

Add-ReferenceRow `
    "Tkinter tutorial" `
    "Website that shows how to use Tkinter" `
    "https://www.geeksforgeeks.org/python-tkinter-tutorial/" `
    "26/4/2024" `
    $false

Add-ReferenceRow `
    "Getting entry value from pressing enter in tkinter " `
    "Stackoverflow discussion about how to get text value from entry fields upon user pressing enter " `
    "https://stackoverflow.com/questions/54846371/pass-a-value-from-a-tkinter-entry-to-a-variable-by-pressing-enter" `
    "27/4/2024" `
    $true

Add-ReferenceRow `
    "Tkinter manual" `
    "Tkinter manual showing some basic functions" `
    "https://www.tutorialspoint.com/python/tk_scrollbar.htm" `
    "29/4/2024" `
    $false

Write-Output ("Table now has " + $t.Rows.Count + " rows")
